$d = $word.ActiveDocument

# --- Paragraph 1: ID marker line ---
$p1 = $d.Paragraphs(1)

# Add a paragraph border (top/left/bottom/right) with 5pt text-to-border spacing
$b = $p1.Format.Borders
$b.DistanceFromTop = 5
$b.DistanceFromLeft = 5
$b.DistanceFromBottom = 5
$b.DistanceFromRight = 5

# Widen the left indent from 120 twips (6pt) to 225 twips (11.25pt)
$p1.Format.LeftIndent = 11.25

# Replace the old topic id text (plus the trailing space run) with the new id,
# collapsing both runs into a single run with no trailing space.
$found = $d.Content.Find.Execute("**ID__AFFARS_pgi_5316_topic_9__ID** ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "**ID__AFFARS_SMC_PGI_5316__ID**", 2)
